$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the first data row (row 2, "H 72"), shifting all subsequent rows up by one.
$ws.Rows.Item(2).Delete()
